$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.230.52"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "2.013.29"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'247.97"
$ws.Range("E5").Value = "  +1.47%  "
$ws.Range("E6").Value = "  -1.78%  "
$ws.Range("D7").Value = "'63.08"
$ws.Range("E7").Value = "  +19.89%  "
$ws.Range("D9").Value = "'59.05"
$ws.Range("E9").Value = "  -3.41%  "
$ws.Range("E10").Value = "  +3.87%  "
$ws.Range("E11").Value = "  +1.72%  "
$ws.Range("E12").Value = "  -0.46%  "
$ws.Range("D13").Value = "'0.955"
$ws.Range("E13").Value = "  +2.69%  "
$ws.Range("D14").Value = "'15.12"
$ws.Range("E14").Value = "  +5.54%  "
$ws.Range("D15").Value = "2.304.99"
$ws.Range("E15").Value = "  -1.12%  "
$ws.Range("E16").Value = "  +2.64%  "
$ws.Range("D17").Value = "'19.59"
$ws.Range("E17").Value = "  +16.66%  "
$ws.Range("D18").Value = "2.014.17"
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("D19").Value = "36.173.19"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").Value = "'72.09"
$ws.Range("E20").Value = "  +1.93%  "
$ws.Range("D21").Value = "0.0₃0860"
$ws.Range("E21").Value = "  +2.26%  "
$ws.Range("D22").Value = "'5.29"
$ws.Range("E22").Value = "  +3.38%  "
$ws.Range("D23").Value = "'235.17"
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").Value = "'2.70"
$ws.Range("E24").Value = "  +22.58%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").Value = "'2.30"
$ws.Range("E26").Value = "  -2.49%  "
$ws.Range("D27").Value = "'9.70"
$ws.Range("E27").Value = "  +6.85%  "
$ws.Range("D28").Value = "'166.21"
$ws.Range("E28").Value = "  +1.76%  "
$ws.Range("D29").Value = "'19.64"
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").Value = "  +0.92%  "
$ws.Range("E31").Value = "  +4.54%  "
$ws.Range("D32").Value = "'1.17"
$ws.Range("E32").Value = "  +0.94%  "
$ws.Range("D33").Value = "'0.100"
$ws.Range("E33").Value = "  +17.04%  "
$ws.Range("D34").Value = "'0.0607"
$ws.Range("E34").Value = "  +3.77%  "
$ws.Range("D35").Value = "'4.51"
$ws.Range("E35").Value = "  +4.40%  "
$ws.Range("D36").Value = "'2.49"
$ws.Range("E36").Value = "  +13.67%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").Value = "'1.81"
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("D39").Value = "'5.79"
$ws.Range("E39").Value = "  +18.68%  "
$ws.Range("E40").Value = "  +2.56%  "
$ws.Range("D41").Value = "'0.0977"
$ws.Range("E41").Value = "  +10.68%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").Value = "'17.26"
$ws.Range("E42").Value = "  +10.91%  "
$ws.Range("E43").Value = "  +2.83%  "
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").Value = "'2.89"
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("E45").Value = "  +3.85%  "
$ws.Range("D46").Value = "'94.79"
$ws.Range("E46").Value = "  +3.13%  "
$ws.Range("D47").Value = "'7.85"
$ws.Range("E47").Value = "  +7.06%  "
$ws.Range("D48").Value = "1.378.75"
$ws.Range("E48").Value = "  +0.41%  "
$ws.Range("D49").Value = "'2.94"
$ws.Range("E49").Value = "  +1.22%  "
$ws.Range("E50").Value = "  +5.30%  "
$ws.Range("D51").Value = "'47.34"
$ws.Range("E51").Value = "  +6.17%  "
